$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (João Rodrigues / Desenho Técnico) - add MEC value
$ws.Range("I2").Value = 2

# Add new row 13 - Josivaldo Ferreira / Programação
$ws.Range("A13").Value = "Josivaldo Ferreira"
$ws.Range("B13").Value = "Programação"
$ws.Range("C13").Value = "Manha"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "S5"
$ws.Range("G13").Value = 1

# Add new row 14 - Pedro Francisco / MTRM (Materia entered before Professor)
$ws.Range("B14").Value = "MTRM"
$ws.Range("A14").Value = "Pedro Francisco"
$ws.Range("C14").Value = "Manha"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "S6"
$ws.Range("F14").Value = "N3"
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 2

# Add new row 15 - Lucas Ferreira / Sistemas digitais (Materia entered before Professor)
$ws.Range("B15").Value = "Sistemas digitais"
$ws.Range("A15").Value = "Lucas Ferreira"
$ws.Range("C15").Value = "Manha"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "S2-N5"
$ws.Range("F15").Value = "N3"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 2

# Add new row 16 - Jorge Aquino / Motores de aplicação (Materia entered before Professor)
$ws.Range("B16").Value = "Motores de aplicação"
$ws.Range("A16").Value = "Jorge Aquino"
$ws.Range("C16").Value = "Manha"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "S3"
$ws.Range("F16").Value = "N6"
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 2

# Update selection to match final cursor location after data entry
$ws.Range("F17").Select()
